# Update the "fraction" values in column C (rows 3 through 402) from 0.01 to 1E-3
# and update the selection to reflect the new active cell/selection range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3:C402").Value = 0.001

$ws.Range("C2:C402").Select()
